$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pollock (row 13): add catch limit values ---
$ws.Range("B13").Value = 11462
$ws.Range("D13").Formula = "=B13/0.2"

# --- Redfish (row 15): rename from "Redfish - 500s" to "Redfish" and add values ---
$ws.Range("A15").Value = "Redfish"
$ws.Range("B15").Value = 526
$ws.Range("D15").Value = 52600

# Give A15 its own distinct style (bold 9pt Arial, left/top aligned, bordered)
# so it no longer shares the generic row-label style used by the other rows.
$ws.Range("A15").Font.Name = "Arial"
$ws.Range("A15").Font.Bold = $true
$ws.Range("A15").Font.Size = 9
$ws.Range("A15").HorizontalAlignment = -4131
$ws.Range("A15").VerticalAlignment = -4160
$ws.Range("A15").Borders.Item(7).LineStyle = 1
$ws.Range("A15").Borders.Item(8).LineStyle = 1
$ws.Range("A15").Borders.Item(9).LineStyle = 1
$ws.Range("A15").Borders.Item(10).LineStyle = 1

# --- White Hake (row 17): add derived TAC formula ---
$ws.Range("D17").Formula = "=B17/0.09"

# --- American Plaice (row 19): add derived TAC formula ---
$ws.Range("D19").Formula = "=B19/0.15"

# --- Witch Flounder (row 25): add derived TAC formula ---
$ws.Range("D25").Formula = "=B25/0.07"

# --- Update selection to match the saved cursor position ---
[void]$ws.Range("D26").Select()
